$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 12,20
$data[0,0] = "ECs"
$data[0,1] = "Sema6d"
$data[0,2] = "Trem2"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 57.38695133333332
$data[0,7] = 172.160854
$data[0,8] = 0.6542464432660997
$data[0,9] = 0.6542464432660998
$data[0,10] = 1
$data[0,11] = 0.3333333333333333
$data[0,12] = 0.03608766666666666
$data[0,13] = 0.108263
$data[0,14] = 0.0004053057704893844
$data[0,15] = 0.0004053057704893844
$data[0,16] = 2.070961170733555
$data[0,17] = 18.638650536602
$data[0,18] = 0.0002651698587779059
$data[0,19] = 0.0002651698587779059

$data[1,0] = "ECs"
$data[1,1] = "Sema6d"
$data[1,2] = "Trem2"
$data[1,3] = "MuSCs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 57.38695133333332
$data[1,7] = 172.160854
$data[1,8] = 0.6542464432660997
$data[1,9] = 0.6542464432660998
$data[1,10] = 1
$data[1,11] = 0.3333333333333333
$data[1,12] = 0.020271
$data[1,13] = 0.060813
$data[1,14] = 0.0002276665141439914
$data[1,15] = 0.0002276665141439913
$data[1,16] = 1.163290890478
$data[1,17] = 10.469618014302
$data[1,18] = 0.0001489500071294975
$data[1,19] = 0.0001489500071294975

$data[2,0] = "ECs"
$data[2,1] = "Sema6d"
$data[2,2] = "Trem2"
$data[2,3] = "Resolving-Mac"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 57.38695133333332
$data[2,7] = 172.160854
$data[2,8] = 0.6542464432660997
$data[2,9] = 0.6542464432660998
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 88.98176833333333
$data[2,13] = 266.945305
$data[2,14] = 0.9993670277153666
$data[2,15] = 0.9993670277153666
$data[2,16] = 5106.39240889894
$data[2,17] = 45957.53168009047
$data[2,18] = 0.6538323234001923
$data[2,19] = 0.6538323234001924

$data[3,0] = "FAPs"
$data[3,1] = "Sema6d"
$data[3,2] = "Trem2"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 14.34807866666667
$data[3,7] = 43.04423600000001
$data[3,8] = 0.16357689713892
$data[3,9] = 0.16357689713892
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.03608766666666666
$data[3,13] = 0.108263
$data[3,14] = 0.0004053057704893844
$data[3,15] = 0.0004053057704893844
$data[3,16] = 0.5177886802297778
$data[3,17] = 4.660098122068001
$data[3,18] = 0.000066298660329152734184236795
$data[3,19] = 0.000066298660329152734184236795

$data[4,0] = "FAPs"
$data[4,1] = "Sema6d"
$data[4,2] = "Trem2"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 14.34807866666667
$data[4,7] = 43.04423600000001
$data[4,8] = 0.16357689713892
$data[4,9] = 0.16357689713892
$data[4,10] = 1
$data[4,11] = 0.3333333333333333
$data[4,12] = 0.020271
$data[4,13] = 0.060813
$data[4,14] = 0.0002276665141439914
$data[4,15] = 0.0002276665141439913
$data[4,16] = 0.2908499026520001
$data[4,17] = 2.617649123868
$data[4,18] = 0.00003724098196610814041429896
$data[4,19] = 0.000037240981966108126861771804

$data[5,0] = "FAPs"
$data[5,1] = "Sema6d"
$data[5,2] = "Trem2"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 14.34807866666667
$data[5,7] = 43.04423600000001
$data[5,8] = 0.16357689713892
$data[5,9] = 0.16357689713892
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 88.98176833333333
$data[5,13] = 266.945305
$data[5,14] = 0.9993670277153666
$data[5,15] = 0.9993670277153666
$data[5,16] = 1276.717411945776
$data[5,17] = 11490.45670751198
$data[5,18] = 0.1634733574966247
$data[5,19] = 0.1634733574966247

$data[6,0] = "MuSCs"
$data[6,1] = "Sema6d"
$data[6,2] = "Trem2"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 14.516389
$data[6,7] = 43.549167
$data[6,8] = 0.1654957381714162
$data[6,9] = 0.1654957381714162
$data[6,10] = 1
$data[6,11] = 0.3333333333333333
$data[6,12] = 0.03608766666666666
$data[6,13] = 0.108263
$data[6,14] = 0.0004053057704893844
$data[6,15] = 0.0004053057704893844
$data[6,16] = 0.5238626074356666
$data[6,17] = 4.714763466920999
$data[6,18] = 0.000067076377672275255066754174
$data[6,19] = 0.000067076377672275255066754174

$data[7,0] = "MuSCs"
$data[7,1] = "Sema6d"
$data[7,2] = "Trem2"
$data[7,3] = "MuSCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 14.516389
$data[7,7] = 43.549167
$data[7,8] = 0.1654957381714162
$data[7,9] = 0.1654957381714162
$data[7,10] = 1
$data[7,11] = 0.3333333333333333
$data[7,12] = 0.020271
$data[7,13] = 0.060813
$data[7,14] = 0.0002276665141439914
$data[7,15] = 0.0002276665141439913
$data[7,16] = 0.294261721419
$data[7,17] = 2.648355492771
$data[7,18] = 0.000037677837815173010683252736
$data[7,19] = 0.00003767783781517299713072558

$data[8,0] = "MuSCs"
$data[8,1] = "Sema6d"
$data[8,2] = "Trem2"
$data[8,3] = "Resolving-Mac"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 14.516389
$data[8,7] = 43.549167
$data[8,8] = 0.1654957381714162
$data[8,9] = 0.1654957381714162
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 88.98176833333333
$data[8,13] = 266.945305
$data[8,14] = 0.9993670277153666
$data[8,15] = 0.9993670277153666
$data[8,16] = 1291.693963034548
$data[8,17] = 11625.24566731094
$data[8,18] = 0.1653909839559287
$data[8,19] = 0.1653909839559287

$data[9,0] = "Resolving-Mac"
$data[9,1] = "Sema6d"
$data[9,2] = "Trem2"
$data[9,3] = "ECs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 1.46316
$data[9,7] = 4.389480000000001
$data[9,8] = 0.01668092142356404
$data[9,9] = 0.01668092142356404
$data[9,10] = 1
$data[9,11] = 0.3333333333333333
$data[9,12] = 0.03608766666666666
$data[9,13] = 0.108263
$data[9,14] = 0.0004053057704893844
$data[9,15] = 0.0004053057704893844
$data[9,16] = 0.05280203036
$data[9,17] = 0.47521827324
$data[9,18] = 0.000006760873710050500595106389
$data[9,19] = 0.000006760873710050500595106389

$data[10,0] = "Resolving-Mac"
$data[10,1] = "Sema6d"
$data[10,2] = "Trem2"
$data[10,3] = "MuSCs"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 1.46316
$data[10,7] = 4.389480000000001
$data[10,8] = 0.01668092142356404
$data[10,9] = 0.01668092142356404
$data[10,10] = 1
$data[10,11] = 0.3333333333333333
$data[10,12] = 0.020271
$data[10,13] = 0.060813
$data[10,14] = 0.0002276665141439914
$data[10,15] = 0.0002276665141439913
$data[10,16] = 0.02965971636
$data[10,17] = 0.26693744724
$data[10,18] = 0.000003797687233212649986473924
$data[10,19] = 0.000003797687233212649986473924

$data[11,0] = "Resolving-Mac"
$data[11,1] = "Sema6d"
$data[11,2] = "Trem2"
$data[11,3] = "Resolving-Mac"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 1.46316
$data[11,7] = 4.389480000000001
$data[11,8] = 0.01668092142356404
$data[11,9] = 0.01668092142356404
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 88.98176833333333
$data[11,13] = 266.945305
$data[11,14] = 0.9993670277153666
$data[11,15] = 0.9993670277153666
$data[11,16] = 130.1945641546
$data[11,17] = 1171.7510773914
$data[11,18] = 0.01667036286262077
$data[11,19] = 0.01667036286262077

$ws.Range("A2:T13").Value = $data
